$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix an existing cell in the second (rows 11-17) puzzle block ---
# F16 was "22c", becomes "27c"
$ws.Range("F16").Value = "27c"

# --- Step 2: build the third puzzle block (rows 22-28), mirroring the
# layout/style of the block at rows 11-17 ---

# Apply the same thin-border style used throughout the sheet to the new block
$ws.Range("A22:G28").Borders.LineStyle = 1

# Row 22: header row
$ws.Range("A22").Value = "a"
$ws.Range("B22").Value = "b"
$ws.Range("C22").Value = "c"
$ws.Range("D22").Value = "d"
$ws.Range("E22").Value = "e"
$ws.Range("F22").Value = "f"
$ws.Range("G22").Value = "on garde"

# Row 23
$ws.Range("A23").Value = 0
$ws.Range("B23").Value = "_"
$ws.Range("C23").Value = "_"
$ws.Range("D23").Value = "_"
$ws.Range("E23").Value = "_"
$ws.Range("F23").Value = "_"
$ws.Range("G23").Value = "a"

# Row 24
$ws.Range("B24").Value = "4a"
$ws.Range("C24").Value = "15a"
$ws.Range("D24").Value = "_"
$ws.Range("E24").Value = "_"
$ws.Range("F24").Value = "_"
$ws.Range("G24").Value = "b"

# Row 25
$ws.Range("C25").Value = "7b"
$ws.Range("D25").Value = "5b"
$ws.Range("E25").Value = "24b"
$ws.Range("F25").Value = "_"
$ws.Range("G25").Value = "d"

# Row 26
$ws.Range("C26").Value = "7d"
$ws.Range("E26").Value = "16d"
$ws.Range("F26").Value = "_"
$ws.Range("G26").Value = "c"

# Row 27
$ws.Range("E27").Value = "16d"
$ws.Range("F27").Value = "12c"
$ws.Range("G27").Value = "f"

# Row 28
$ws.Range("E28").Value = "14f"
$ws.Range("G28").Value = "e"

# --- Step 3: trailing notes below the new block (rows 31-34), column G only,
# no border styling (matches plain unstyled cells elsewhere on the sheet) ---
$ws.Range("G31").Value = "a-b-d-c-f-e"
$ws.Range("G32").Value = "ou"
$ws.Range("G33").Value = "a-b-c-f-e"
$ws.Range("G34").Value = "car meme valeur"

# --- Step 4: move the active selection, as reflected in the saved file ---
$ws.Range("K32").Select()
